$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 5947.1
$ws.Range("I28").Value = 5496.778
$ws.Range("J28").Value = 10000
$ws.Range("K28").Value = 5496.778
$ws.Range("L28").Value = 10000
$ws.Range("M28").Value = -5011.778
$ws.Range("N28").Value = -10970

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 22732832
$ws.Range("I62").Value = 27782462
$ws.Range("K62").Value = 27782462
$ws.Range("M62").Value = -27781838

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 22732832
$ws.Range("I65").Value = 27782462
$ws.Range("K65").Value = 138912310
$ws.Range("M65").Value = -138909190

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 951.63635
$ws.Range("I96").Value = 674.2
$ws.Range("J96").Value = 1182.8334
$ws.Range("K96").Value = 2022.6
$ws.Range("L96").Value = 3548.5002
$ws.Range("M96").Value = -649.6000000000001
$ws.Range("N96").Value = -6294.5002

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 9233.972
$ws.Range("I98").Value = 1725.6818
$ws.Range("J98").Value = 21940.309
$ws.Range("K98").Value = 1725.6818
$ws.Range("L98").Value = 21940.309
$ws.Range("M98").Value = -227.6818000000001
$ws.Range("N98").Value = -24936.309

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 12998431
$ws.Range("I106").Value = 18195810
$ws.Range("K106").Value = 18195810
$ws.Range("M106").Value = -18195179

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 9233.972
$ws.Range("I122").Value = 1725.6818
$ws.Range("J122").Value = 21940.309
$ws.Range("K122").Value = 5177.0454
$ws.Range("L122").Value = 65820.927
$ws.Range("M122").Value = -2727.0454
$ws.Range("N122").Value = -70720.927

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1332.4
$ws.Range("I137").Value = 1141.1
$ws.Range("K137").Value = 3423.3
$ws.Range("M137").Value = -873.2999999999997

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3560.507
$ws.Range("I138").Value = 4115.2
$ws.Range("J138").Value = 3259.0435
$ws.Range("K138").Value = 12345.6
$ws.Range("L138").Value = 9777.130500000001
$ws.Range("M138").Value = -7205.599999999999
$ws.Range("N138").Value = -20057.1305

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4996.61
$ws.Range("I32").Value = 3684
$ws.Range("J32").Value = 12653.5
$ws.Range("K32").Value = 3684
$ws.Range("L32").Value = 12653.5
$ws.Range("M32").Value = -3397
$ws.Range("N32").Value = -13227.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2629.1428
$ws.Range("I61").Value = 2280
$ws.Range("J61").Value = 3502
$ws.Range("K61").Value = 2280
$ws.Range("L61").Value = 3502
$ws.Range("M61").Value = -2068
$ws.Range("N61").Value = -3926

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H95").Value = 47598
$ws.Range("J95").Value = 47598
$ws.Range("L95").Value = 47598
$ws.Range("N95").Value = -53090

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 4638.1304
$ws.Range("I102").Value = 4638.1304
$ws.Range("K102").Value = 4638.1304
$ws.Range("M102").Value = -3016.1304

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4330.8857
$ws.Range("I132").Value = 4289.909
$ws.Range("J132").Value = 5007
$ws.Range("K132").Value = 12869.727
$ws.Range("L132").Value = 15021
$ws.Range("M132").Value = -10339.727
$ws.Range("N132").Value = -20081

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2629.1428
$ws.Range("I136").Value = 2280
$ws.Range("J136").Value = 3502
$ws.Range("K136").Value = 6840
$ws.Range("L136").Value = 10506
$ws.Range("M136").Value = -4290
$ws.Range("N136").Value = -15606

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3956.8823
$ws.Range("I105").Value = 2790.889
$ws.Range("K105").Value = 2790.889
$ws.Range("M105").Value = -1043.889

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6343.423
$ws.Range("I134").Value = 5304.067
$ws.Range("J134").Value = 7760.727
$ws.Range("K134").Value = 15912.201
$ws.Range("L134").Value = 23282.181
$ws.Range("M134").Value = -13377.201
$ws.Range("N134").Value = -28352.181

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 43213.2
$ws.Range("I31").Value = 1865.1818
$ws.Range("J31").Value = 75700.92999999999
$ws.Range("K31").Value = 1865.1818
$ws.Range("L31").Value = 75700.92999999999
$ws.Range("M31").Value = -1570.1818
$ws.Range("N31").Value = -76290.92999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 43213.2
$ws.Range("I34").Value = 1865.1818
$ws.Range("J34").Value = 75700.92999999999
$ws.Range("K34").Value = 1865.1818
$ws.Range("L34").Value = 75700.92999999999
$ws.Range("M34").Value = -1663.1818
$ws.Range("N34").Value = -76104.92999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2546.1428
$ws.Range("I58").Value = 2194
$ws.Range("J58").Value = 3015.6667
$ws.Range("K58").Value = 2194
$ws.Range("L58").Value = 3015.6667
$ws.Range("M58").Value = -1991
$ws.Range("N58").Value = -3421.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3373.258
$ws.Range("I122").Value = 3143.75
$ws.Range("K122").Value = 9431.25
$ws.Range("M122").Value = -6981.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 26631.95
$ws.Range("I134").Value = 38510.742
$ws.Range("J134").Value = 1960.6154
$ws.Range("K134").Value = 115532.226
$ws.Range("L134").Value = 5881.8462
$ws.Range("M134").Value = -112997.226
$ws.Range("N134").Value = -10951.8462

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2546.1428
$ws.Range("I136").Value = 2194
$ws.Range("J136").Value = 3015.6667
$ws.Range("K136").Value = 6582
$ws.Range("L136").Value = 9047.000100000001
$ws.Range("M136").Value = -4032
$ws.Range("N136").Value = -14147.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 344467.5
$ws.Range("I141").Value = 46000
$ws.Range("J141").Value = 523548
$ws.Range("K141").Value = 46000
$ws.Range("L141").Value = 523548
$ws.Range("M141").Value = -40820
$ws.Range("N141").Value = -533908

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 581
$ws.Range("I2").Value = 299
$ws.Range("K2").Value = 1794
$ws.Range("M2").Value = -1681

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 9317722
$ws.Range("I4").Value = 10073878
$ws.Range("K4").Value = 30221634
$ws.Range("M4").Value = -30221522

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 184.6
$ws.Range("I7").Value = 99.5
$ws.Range("J7").Value = 312.25
$ws.Range("K7").Value = 298.5
$ws.Range("L7").Value = 936.75
$ws.Range("M7").Value = -186.5
$ws.Range("N7").Value = -1160.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 140.2
$ws.Range("I14").Value = 140.2
$ws.Range("K14").Value = 420.6
$ws.Range("M14").Value = -247.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 224.36363
$ws.Range("I38").Value = 144.8
$ws.Range("J38").Value = 290.66666
$ws.Range("K38").Value = 434.4
$ws.Range("L38").Value = 871.9999799999999
$ws.Range("M38").Value = -87.40000000000003
$ws.Range("N38").Value = -1565.99998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 66
$ws.Range("I40").Value = 54.2
$ws.Range("K40").Value = 216.8
$ws.Range("M40").Value = -147.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 752
$ws.Range("I92").Value = 709.6
$ws.Range("J92").Value = 794.4
$ws.Range("K92").Value = 2128.8
$ws.Range("L92").Value = 2383.2
$ws.Range("M92").Value = -880.8000000000002
$ws.Range("N92").Value = -4879.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1685.8182
$ws.Range("I122").Value = 498.75
$ws.Range("J122").Value = 2364.1428
$ws.Range("K122").Value = 4488.75
$ws.Range("L122").Value = 21277.2852
$ws.Range("M122").Value = -2038.75
$ws.Range("N122").Value = -26177.2852

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 16199.25
$ws.Range("J92").Value = 16199.25
$ws.Range("L92").Value = 16199.25
$ws.Range("N92").Value = -19943.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 23669.25
$ws.Range("I122").Value = 34028.395
$ws.Range("K122").Value = 102085.185
$ws.Range("M122").Value = -99635.185

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5930.1304
$ws.Range("I126").Value = 4883.9375
$ws.Range("J126").Value = 8321.429
$ws.Range("K126").Value = 14651.8125
$ws.Range("L126").Value = 24964.287
$ws.Range("M126").Value = -12181.8125
$ws.Range("N126").Value = -29904.287

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5518.6
$ws.Range("I7").Value = 5518.6
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 5518.6
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -5406.6
$ws.Range("N7").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2959.7856
$ws.Range("I40").Value = 2735.818
$ws.Range("K40").Value = 2735.818
$ws.Range("M40").Value = -2599.818

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1008.5909
$ws.Range("I55").Value = 216.625
$ws.Range("K55").Value = 216.625
$ws.Range("M55").Value = -43.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 5518.6
$ws.Range("I126").Value = 5518.6
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 16555.8
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -14085.8
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3960.9592
$ws.Range("I132").Value = 3631.3333
$ws.Range("K132").Value = 10893.9999
$ws.Range("M132").Value = -8363.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 100004.8
$ws.Range("J140").Value = 100004.8
$ws.Range("L140").Value = 100004.8
$ws.Range("N140").Value = -110364.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 5474.5
$ws.Range("J15").Value = 5474.5
$ws.Range("L15").Value = 5474.5
$ws.Range("N15").Value = -6050.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2310.2104
$ws.Range("I126").Value = 2046.0714
$ws.Range("K126").Value = 6138.2142
$ws.Range("M126").Value = -3668.2142

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2141.12
$ws.Range("I132").Value = 1938.6666
$ws.Range("K132").Value = 5815.9998
$ws.Range("M132").Value = -3285.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5058.2925
$ws.Range("I136").Value = 4027.4546
$ws.Range("J136").Value = 9310.5
$ws.Range("K136").Value = 12082.3638
$ws.Range("L136").Value = 27931.5
$ws.Range("M136").Value = -9532.363799999999
$ws.Range("N136").Value = -33031.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H140").Value = 83851.336
$ws.Range("J140").Value = 83851.336
$ws.Range("L140").Value = 83851.336
$ws.Range("N140").Value = -94211.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 99333.31
$ws.Range("I141").Value = 85500
$ws.Range("J141").Value = 100255.53
$ws.Range("K141").Value = 85500
$ws.Range("L141").Value = 100255.53
$ws.Range("M141").Value = -80320
$ws.Range("N141").Value = -110615.53
